$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 6 (the child "Corene/Myra" record) - this shifts all subsequent
# child rows, the school row, cost row and time row up by one, and updates the
# sheet dimension from A1:H16 to A1:H15.
$ws.Rows(6).Delete()

# Update the number of children (the route was recalculated for 7 children).
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = '7'

# The whole route (stop order, coordinates, pickup times, parent contacts and
# remaining distances/times) was recalculated by the path-finding routine, so
# write out the freshly computed values for every child row plus the school,
# cost and total time rows. Force text formatting first so that values such as
# "30.0" or "7" are kept as text instead of being converted to numbers, matching
# the original inline-string layout of the sheet.

# Row 6
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = '0'
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = '5'
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'Patti  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = 'Lavenia  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-6.55,4.12'
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = 'Jennell(mother): 0503029941'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '7:00:00'
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = '30.0'

# Row 7
$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = '1'
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = '9'
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'Letha  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = 'Stephenie  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '-8.77,7.51'
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = 'Sibyl(mother): 0567328221'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '7:05:00'
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = '25.0'

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = '2'
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = '7'
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'Wyatt  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = 'Willette  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-7.11,9.53'
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = 'Antionette(father): 0557331799'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '7:09:00'
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = '21.0'

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = '3'
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = '17'
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'Britta  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = 'Jamel  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-3.64,7.49'
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = 'Albertine(father): 0574981040'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '7:14:00'
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = '16.0'

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = '4'
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = '12'
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'Frankie  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = 'Flavia  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '1.45,4.78'
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = 'Cyrus(mother): 0522363358'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '7:21:00'
$ws.Range("H10").NumberFormat = "@"
$ws.Range("H10").Value = '9.0'

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = '5'
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = '4'
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'Francisca  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = 'Stevie  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-0.28,2.19'
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = 'Bernardine(mother): 0561339273'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '7:25:00'
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = '5.0'

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = '6'
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = '13'
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'Fay  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = 'Emilee  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '0.33,2.02'
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = 'Sheri(mother): 0516797453'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '7:26:00'
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = '4.0'

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = 'school'
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = '3'
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'Ironiah'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = 'mySchool'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '0,0'
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = 'Shir(secretary): 0523345098'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '7:30:00'

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = 'cost'
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = '25'

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = 'time'
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = '30.0'
